$wb = $excel.ActiveWorkbook

# Overview sheet: update "Latest HO Xliff Generate Date" for rows 4-7
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("G4:G7").Value = "2016-09-01 18:36:39"

# zh-cn sheet: Priority low -> ht, Latest Handoff Datetime updated
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("E4:E7").Value = "ht"
$ws2.Range("H4:H7").Value = "2016-09-01 18:36:34"

# de-de sheet: Priority low -> ht
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("E4:E7").Value = "ht"
